# Edit slide 6 ("얻어 낼 수 있는 점?") body placeholder:
#  - Rewrite part of bullet 2's text (shorten/replace the parenthetical
#    explanation with new wording around "weight").
#  - Append a brand-new bullet "4. ..." after the existing bullet 3.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# --- Bullet 2 ------------------------------------------------------------
# Current runs in paragraph 2:
#   1 "2. "
#   2 "스마트폰과 와이파이의 우선연결과정"
#   3 "("
#   4 "즉 성능의 좋은 휴대전화의 경우 먼저 연결되는 과정"
#   5 ")"
#   6 "에 어떤 방식이 사용되는지를 "
#   7 "weight"
#   8 "를 이용해 알아볼 수 있다"
#   9 "."
$para2 = $tr.Paragraphs(2, 1)

# Remove runs 3-6 (the "(...)  에 어떤 방식이 사용되는지를 " runs); deleting a
# single run at a time collapses it out of the XML, so always target
# position 3 since everything shifts down after each delete.
$para2.Runs(3, 1).Text = ""
$para2.Runs(3, 1).Text = ""
$para2.Runs(3, 1).Text = ""
$para2.Runs(3, 1).Text = ""

# Replace run 2's text.
$para2.Runs(2, 1).Text = "스마트폰에 잡힌 여러 와이파이들을 중 와이파이 강도가 쌘 와이파이에 "

# Replace what used to be run 8 (now run 4, after "weight").
$para2.Runs(4, 1).Text = "를 줘서 스마트폰이 원활한 데이터이용을 할 수 있다"

# --- New bullet 4 ----------------------------------------------------------
$para3 = $tr.Paragraphs(3, 1)
$null = $para3.InsertAfter("`r4. 코드를 개선한 경우를 통해 좀 더 휴대전화가 와이파이에 최적으로 접근하고 한 와이파이에 휴대전화들이 뭉쳐 느려지는 것을 보안할 수 있다.")

# Split the new paragraph's single run into the same run boundaries as the
# reference edit: "4. " | "코드를...뭉쳐 " | "느려지는" | " 것을...있다" | "."
$para4 = $tr.Paragraphs(4, 1)
$full  = $para4.Runs(1, 1).Text
$t1 = "4. "
$t2 = "코드를 개선한 경우를 통해 좀 더 휴대전화가 와이파이에 최적으로 접근하고 한 와이파이에 휴대전화들이 뭉쳐 "
$t3 = "느려지는"
$t4 = " 것을 보안할 수 있다"
$t5 = "."

$o1 = 1
$o2 = $o1 + $t1.Length
$o3 = $o2 + $t2.Length
$o4 = $o3 + $t3.Length
$o5 = $o4 + $t4.Length

# Re-assigning a sub-range's Text to itself forces the engine to split the
# run boundary there (without changing the actual characters).
$c1 = $para4.Characters($o1, $t1.Length); $c1.Text = $c1.Text
$c2 = $para4.Characters($o2, $t2.Length); $c2.Text = $c2.Text
$c3 = $para4.Characters($o3, $t3.Length); $c3.Text = $c3.Text
$c4 = $para4.Characters($o4, $t4.Length); $c4.Text = $c4.Text
$c5 = $para4.Characters($o5, $t5.Length); $c5.Text = $c5.Text
